# Auto-generated edit script: add "Pleine Lune" (Full Moon) Rahu gifts
# rows 37-41 on sheet "Dons d'origine" (Worksheets.Item(2))
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# ---- Row 37 ----
$ws.Range("C37").Value = @'
Instinct de Tueur
'@
$ws.Range("E37").Value = @'
Le Rahu est un tueur forgé par la Lune. Cette Facette met en avant ses instincts mortels
'@
$ws.Range("F37").Value = @'
E
'@
$ws.Range("G37").Value = @'
n/a
'@
$ws.Range("H37").Value = @'
Reflexe
'@
$ws.Range("I37").Value = @'
1 scène
'@
$ws.Range("J37").Value = @'
La Rahu bénéficie du 8-relance sur tous ses jets de Bagarre et Mélée pendant la durée de la scène
'@
$ws.Range("K37").Value = @'
n/a
'@
$ws.Range("L37").Value = @'
n/a
'@
$ws.Range("M37").Value = @'
n/a
'@
$ws.Range("N37").Value = @'
n/a
'@

# ---- Row 38 ----
$ws.Range("C38").Value = @'
Peau de Guerrier
'@
$ws.Range("E38").Value = @'
La bataille est une épreuve. La Rahu doit être suffisament endurante pour encaisser le pire que son adversaire peut lui donner et donner le coup de grâce en retour
'@
$ws.Range("F38").Value = @'
n/a
'@
$ws.Range("G38").Value = @'
n/a
'@
$ws.Range("H38").Value = @'
n/a
'@
$ws.Range("I38").Value = @'
Permanent
'@
$ws.Range("J38").Value = @'
La Rahu ajoute sa Renommée Pureté à sa Santé, gagnant ainsi une case Santé à chaque fois qu'elle augmente sa Pureté. Cette augmentation est permanente.
'@
$ws.Range("K38").Value = @'
n/a
'@
$ws.Range("L38").Value = @'
n/a
'@
$ws.Range("M38").Value = @'
n/a
'@
$ws.Range("N38").Value = @'
n/a
'@

# ---- Row 39 ----
$ws.Range("C39").Value = @'
Chasseur aux Mains Couvertes de Sang
'@
$ws.Range("E39").Value = @'
Une Rahu en chasse est une chose terrifiante à observer, une figure couverte de sang entièrement débouée au meurtre. Malheure à quiconque se met sur son chemin.
Cette facette ne peut être activée que quand la Rahu acquiers l'Etat Siskur-Dah.
'@
$ws.Range("E39").Characters(231, 11).Font.Italic = $true
$ws.Range("F39").Value = @'
E
'@
$ws.Range("G39").Value = @'
n/a
'@
$ws.Range("H39").Value = @'
Instantanée
'@
$ws.Range("I39").Value = @'
Toute la durée de la Siskur-Dah
'@
$ws.Range("J39").Value = @'
Pendant la durée de la Facette la Rahu ajoute sa Renommée Pureté à ses jets d'attaque contre tout obstacle à la Siskur-Dah. Que la cible soit un allié de la cible qui cherche à tuer la Rahu ou que ce soit un petit bureaucrate dont le métier fait obstruction la Rahu gagne le bonus. Cette Facette ne donne pas de bonus contre la cible de la Siskur-Dah elle même.
'@
$ws.Range("J39").Characters(113, 12).Font.Italic = $true
$ws.Range("J39").Characters(341, 10).Font.Italic = $true
$ws.Range("K39").Value = @'
n/a
'@
$ws.Range("L39").Value = @'
n/a
'@
$ws.Range("M39").Value = @'
n/a
'@
$ws.Range("N39").Value = @'
n/a
'@

# ---- Row 40 ----
$ws.Range("C40").Value = @'
Boucherie
'@
$ws.Range("E40").Value = @'
Grâce à cette Facette la Rahu excelle quand elle affronte plusieurs ennemis. C'est une opportunité pour peindre les murs en rouge.
'@
$ws.Range("F40").Value = @'
EE
'@
$ws.Range("G40").Value = @'
Astuce + Bagarre + Pureté
'@
$ws.Range("H40").Value = @'
Reflexe
'@
$ws.Range("I40").Value = @'
1 tour par succès
'@
$ws.Range("J40").Value = @'
n/a
'@
$ws.Range("K40").Value = @'
La Rahu assoifée de sang se surménage. Pendant un tour sa Défense tombe à zéro et elle ne peux pas entreprendre d'actions qui la feraient abandonner sa Défense.
'@
$ws.Range("L40").Value = @'
Aucun effet
'@
$ws.Range("M40").Value = @'
Pour la durée de la Facette et tant que la Rahu combat plus d'un ennemi, à chaque fois que la Rahu attaque un ennemi avec une attaque sans armes (griffes et cros inclus), ou qu'elle est frappée par un adversaire, elle peut infliger l'un des Etats suivants même si l'attaque ne fait pas de dégâts : "Arm Wrack", "Blinded", "Deafened", "Knocked Down", "Leg Wrack". A chaque fois que la Rahu met un adversaire hors combat en le tuant ou en l'incapacitant la durée de cette Facette est prolongée d'un tour.
'@
$ws.Range("N40").Value = @'
A la première application la Rahu applique 3 des Etats.
'@

# ---- Row 41 ----
$ws.Range("C41").Value = @'
Spasme Ecarlate
'@
$ws.Range("E41").Value = @'
Le Rahu laisse l'entière et pimitive délire de la Lune l'emplir avec de la rage transformatrice. Sa chair se plie, tords et spasme en une forme primale de destruction.
'@
$ws.Range("F41").Value = @'
1 Essence par tour
'@
$ws.Range("G41").Value = @'
Vigueur + Survie + Pureté
'@
$ws.Range("H41").Value = @'
Instantanée
'@
$ws.Range("I41").Value = @'
n/a
'@
$ws.Range("J41").Value = @'
n/a
'@
$ws.Range("K41").Value = @'
La Rahu perds le contrôle de la puissance brute de la lune qui déferle en elle. Ceci compte comme un "breaking point" vers l'Esprit et la Rahu se métamorphoe involontairement dans une autre forme.
'@
$ws.Range("L41").Value = @'
La Facette échoue
'@
$ws.Range("M41").Value = @'
Chaque réussite permet à la Rahu d'ajouter deux points de Force ou Vigueur, d'ajouter un point d'armure générale, ou d'augmenter les dégâts létaux de ses armes naturelles par +1. Ces bénéfices s'accumulent à ceux octroyés par d'autres Facettes et formes. Le Spasme dure tant que la Rahu l'entretiens au prix d'un point d'Essence par tout. Si la Rahu utilise cette Facette sous la forme Hishu ou Urhan elle provoque le Delire comme si elle était sous la forme Dalu.
'@
$ws.Range("N41").Value = @'
La Rahu régènere aussi un point de dégâts létaux à chaque tour pour toute la durée du Spasme.
'@

# ---- Row heights (approximate autofit results from the authored workbook) ----
$ws.Rows.Item(37).RowHeight = 30
$ws.Rows.Item(38).RowHeight = 60
$ws.Rows.Item(39).RowHeight = 105
$ws.Rows.Item(40).RowHeight = 150
$ws.Rows.Item(41).RowHeight = 150

# ---- Selection / view state ----
$ws.Activate()
$ws.Range("C42").Select()
